# Refresh the cryptocurrency Price (D) and Volume/1h change (E) columns
# with the latest scraped values, preserving each cell's existing style.
$ws = $excel.ActiveWorkbook.ActiveSheet

function Set-TextValue {
    # Assigns $Text to $Cell as literal text, even when it looks like a number
    # (e.g. "566.06"), by briefly switching to a Text number format so Excel
    # does not silently convert it to a floating-point value.
    param($Cell, $Text)
    $originalStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $originalStyle
}

$ws.Range('D2').Value = '59.450.18'
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('D3').Value = '2.584.64'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  +0.20%  '
Set-TextValue $ws.Range('D5') '566.06'
$ws.Range('E5').Value = '  +4.44%  '
Set-TextValue $ws.Range('D6') '142.78'
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').Value = '  +0.04%  '
Set-TextValue $ws.Range('D8') '0.595'
$ws.Range('E8').Value = '  +2.68%  '
$ws.Range('D9').Value = '2.588.20'
$ws.Range('E9').Value = '  +0.41%  '
Set-TextValue $ws.Range('D10') '6.64'
$ws.Range('E10').Value = '  -1.72%  '
$ws.Range('E11').Value = '  +2.54%  '
Set-TextValue $ws.Range('D12') '0.152'
$ws.Range('E12').Value = '  +9.75%  '
Set-TextValue $ws.Range('D13') '0.340'
$ws.Range('E13').Value = '  +2.58%  '
$ws.Range('D14').Value = '3.048.71'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('D15').Value = '59.511.44'
$ws.Range('E15').Value = '  +2.02%  '
Set-TextValue $ws.Range('D16') '21.82'
$ws.Range('E16').Value = '  +6.20%  '
$ws.Range('E17').Value = '  +3.55%  '
$ws.Range('D18').Value = '2.597.58'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('E19').Value = '  +1.13%  '
Set-TextValue $ws.Range('D20') '336.39'
$ws.Range('E20').Value = '  +0.65%  '
Set-TextValue $ws.Range('D21') '10.18'
$ws.Range('E21').Value = '  +1.68%  '
Set-TextValue $ws.Range('D22') '6.22'
$ws.Range('E22').Value = '  +2.13%  '
Set-TextValue $ws.Range('D24') '65.03'
$ws.Range('E24').Value = '  -2.05%  '
Set-TextValue $ws.Range('D25') '0.444'
$ws.Range('E25').Value = '  +5.53%  '
Set-TextValue $ws.Range('D26') '0.999'
$ws.Range('E26').Value = '  +0.04%  '
Set-TextValue $ws.Range('D27') '0.161'
$ws.Range('E27').Value = '  +1.91%  '
Set-TextValue $ws.Range('D28') '7.24'
$ws.Range('E28').Value = '  +2.78%  '
$ws.Range('E29').Value = '  +6.49%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  +3.27%  '
Set-TextValue $ws.Range('D32') '160.27'
$ws.Range('E32').Value = '  +4.61%  '
Set-TextValue $ws.Range('D33') '6.02'
$ws.Range('E33').Value = '  +1.04%  '
Set-TextValue $ws.Range('D34') '18.93'
$ws.Range('E34').Value = '  +0.10%  '
Set-TextValue $ws.Range('D35') '4.01'
$ws.Range('E35').Value = '  +3.12%  '
Set-TextValue $ws.Range('D36') '0.885'
$ws.Range('E36').Value = '  +8.40%  '
Set-TextValue $ws.Range('D37') '0.875'
$ws.Range('E37').Value = '  +2.92%  '
$ws.Range('E38').Value = '  +3.65%  '
Set-TextValue $ws.Range('D39') '37.02'
$ws.Range('E39').Value = '  +0.02%  '
Set-TextValue $ws.Range('D40') '1.48'
$ws.Range('E40').Value = '  +4.82%  '
Set-TextValue $ws.Range('D41') '292.52'
$ws.Range('E41').Value = '  +4.93%  '
$ws.Range('E42').Value = '  +1.30%  '
Set-TextValue $ws.Range('D43') '0.998'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  +3.54%  '
Set-TextValue $ws.Range('D45') '0.590'
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('E46').Value = '  +0.28%  '
$ws.Range('E47').Value = '  +0.95%  '
Set-TextValue $ws.Range('D48') '18.97'
$ws.Range('E48').Value = '  +2.89%  '
Set-TextValue $ws.Range('D49') '124.50'
$ws.Range('E49').Value = '  +14.75%  '
Set-TextValue $ws.Range('D50') '0.0231'
$ws.Range('E50').Value = '  +2.48%  '
$ws.Range('D51').Value = '1.928.49'
$ws.Range('E51').Value = '  +1.49%  '
